$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4769.6
$ws.Range("I29").Value = 303.8
$ws.Range("J29").Value = 9235.4
$ws.Range("K29").Value = 911.4000000000001
$ws.Range("L29").Value = 27706.2
$ws.Range("M29").Value = -630.4000000000001
$ws.Range("N29").Value = -28268.2

$ws.Range("H31").Value = 1798.3334
$ws.Range("I31").Value = 997.5
$ws.Range("K31").Value = 2992.5
$ws.Range("M31").Value = -2762.5

$ws.Range("H33").Value = 531.38464
$ws.Range("I33").Value = 586.375
$ws.Range("K33").Value = 586.375
$ws.Range("M33").Value = -357.375

$ws.Range("H49").Value = 2437.875
$ws.Range("I49").Value = 616.6667
$ws.Range("J49").Value = 3530.6
$ws.Range("K49").Value = 1850.0001
$ws.Range("L49").Value = 10591.8
$ws.Range("M49").Value = -1714.0001
$ws.Range("N49").Value = -10863.8

$ws.Range("H51").Value = 8345.385
$ws.Range("J51").Value = 8839.200000000001
$ws.Range("L51").Value = 8839.200000000001
$ws.Range("N51").Value = -9807.200000000001

$ws.Range("H86").Value = 4559.2
$ws.Range("I86").Value = 3401.5
$ws.Range("J86").Value = 4848.625
$ws.Range("K86").Value = 3401.5
$ws.Range("L86").Value = 4848.625
$ws.Range("M86").Value = -2278.5
$ws.Range("N86").Value = -7094.625

$ws.Range("H89").Value = 4559.2
$ws.Range("I89").Value = 3401.5
$ws.Range("J89").Value = 4848.625
$ws.Range("K89").Value = 17007.5
$ws.Range("L89").Value = 24243.125
$ws.Range("M89").Value = -11391.5
$ws.Range("N89").Value = -35475.125

$ws.Range("H112").Value = 3057.4285
$ws.Range("J112").Value = 3800.6
$ws.Range("L112").Value = 11401.8
$ws.Range("N112").Value = -13617.8

$ws.Range("H113").Value = 6945.8335
$ws.Range("J113").Value = 7119.706
$ws.Range("L113").Value = 7119.706
$ws.Range("N113").Value = -13627.706

$ws.Range("H127").Value = 1378.75
$ws.Range("I127").Value = 506.42856
$ws.Range("J127").Value = 2600
$ws.Range("K127").Value = 1519.28568
$ws.Range("L127").Value = 7800
$ws.Range("M127").Value = 3440.71432
$ws.Range("N127").Value = -17720

$ws.Range("H132").Value = 1574.4082
$ws.Range("I132").Value = 1280.8695
$ws.Range("K132").Value = 3842.6085
$ws.Range("M132").Value = -1312.6085

$ws.Range("H133").Value = 60814
$ws.Range("J133").Value = 60814
$ws.Range("L133").Value = 60814
$ws.Range("N133").Value = -70934

$ws.Range("H134").Value = 65251.555
$ws.Range("J134").Value = 65251.555
$ws.Range("L134").Value = 65251.555
$ws.Range("N134").Value = -75391.55499999999

$ws.Range("H137").Value = 26318934
$ws.Range("J137").Value = 3425.1853
$ws.Range("L137").Value = 10275.5559
$ws.Range("N137").Value = -15375.5559

$ws.Range("H138").Value = 3033.0195
$ws.Range("I138").Value = 1612.6552
$ws.Range("J138").Value = 4905.3184
$ws.Range("K138").Value = 4837.9656
$ws.Range("L138").Value = 14715.9552
$ws.Range("M138").Value = 302.0344000000005
$ws.Range("N138").Value = -24995.9552

$ws.Range("H141").Value = 1295.5
$ws.Range("I141").Value = 1095.1364
$ws.Range("J141").Value = 3499.5
$ws.Range("K141").Value = 3285.4092
$ws.Range("L141").Value = 10498.5
$ws.Range("M141").Value = 1894.5908
$ws.Range("N141").Value = -20858.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 419.42856
$ws.Range("I5").Value = 398
$ws.Range("J5").Value = 423
$ws.Range("K5").Value = 398
$ws.Range("L5").Value = 423
$ws.Range("M5").Value = -286
$ws.Range("N5").Value = -647

$ws.Range("H45").Value = 3310.2
$ws.Range("I45").Value = 1808
$ws.Range("K45").Value = 1808
$ws.Range("M45").Value = -1431

$ws.Range("H60").Value = 42000
$ws.Range("I60").Value = 42000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 42000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -41267
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 9410.837
$ws.Range("I61").Value = 6375.816
$ws.Range("K61").Value = 6375.816
$ws.Range("M61").Value = -6163.816

$ws.Range("H74").Value = 25644242
$ws.Range("I74").Value = 37039464
$ws.Range("K74").Value = 37039464
$ws.Range("M74").Value = -37038590

$ws.Range("H77").Value = 25644242
$ws.Range("I77").Value = 37039464
$ws.Range("K77").Value = 185197320
$ws.Range("M77").Value = -185192952

$ws.Range("H97").Value = 1067.25
$ws.Range("I97").Value = 1071.7333
$ws.Range("K97").Value = 1071.7333
$ws.Range("M97").Value = -575.7333000000001

$ws.Range("H124").Value = 38528.285
$ws.Range("J124").Value = 38528.285
$ws.Range("L124").Value = 38528.285
$ws.Range("N124").Value = -48348.285

$ws.Range("H132").Value = 2508.0425
$ws.Range("I132").Value = 1999.425
$ws.Range("K132").Value = 5998.275
$ws.Range("M132").Value = -3468.275

$ws.Range("H135").Value = 73773.60000000001
$ws.Range("J135").Value = 73773.60000000001
$ws.Range("L135").Value = 73773.60000000001
$ws.Range("N135").Value = -83913.60000000001

$ws.Range("H136").Value = 9410.837
$ws.Range("I136").Value = 6375.816
$ws.Range("K136").Value = 19127.448
$ws.Range("M136").Value = -16577.448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 419.42856
$ws.Range("I4").Value = 398
$ws.Range("J4").Value = 423
$ws.Range("K4").Value = 398
$ws.Range("L4").Value = 423
$ws.Range("M4").Value = -283
$ws.Range("N4").Value = -653

$ws.Range("H22").Value = 841.25
$ws.Range("I22").Value = 788.3333
$ws.Range("K22").Value = 788.3333
$ws.Range("M22").Value = -615.3333

$ws.Range("H105").Value = 9109
$ws.Range("I105").Value = 8187.871
$ws.Range("J105").Value = 13868.167
$ws.Range("K105").Value = 8187.871
$ws.Range("L105").Value = 13868.167
$ws.Range("M105").Value = -6440.871
$ws.Range("N105").Value = -17362.167

$ws.Range("H107").Value = 2254.4546
$ws.Range("I107").Value = 915
$ws.Range("K107").Value = 915
$ws.Range("M107").Value = 1005

$ws.Range("H132").Value = 66330.664
$ws.Range("J132").Value = 66330.664
$ws.Range("L132").Value = 66330.664
$ws.Range("N132").Value = -76450.664

$ws.Range("H134").Value = 2279.0688
$ws.Range("I134").Value = 1118.7391
$ws.Range("K134").Value = 3356.2173
$ws.Range("M134").Value = -821.2173000000003

$ws.Range("H135").Value = 45262.5
$ws.Range("J135").Value = 45262.5
$ws.Range("L135").Value = 45262.5
$ws.Range("N135").Value = -55402.5

$ws.Range("H138").Value = 65494
$ws.Range("J138").Value = 65494
$ws.Range("L138").Value = 65494
$ws.Range("N138").Value = -75774

$ws.Range("H140").Value = 58712.2
$ws.Range("J140").Value = 58712.2
$ws.Range("L140").Value = 58712.2
$ws.Range("N140").Value = -69072.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 49437.78
$ws.Range("I31").Value = 1478.4615
$ws.Range("J31").Value = 111784.9
$ws.Range("K31").Value = 1478.4615
$ws.Range("L31").Value = 111784.9
$ws.Range("M31").Value = -1183.4615
$ws.Range("N31").Value = -112374.9

$ws.Range("H34").Value = 49437.78
$ws.Range("I34").Value = 1478.4615
$ws.Range("J34").Value = 111784.9
$ws.Range("K34").Value = 1478.4615
$ws.Range("L34").Value = 111784.9
$ws.Range("M34").Value = -1276.4615
$ws.Range("N34").Value = -112188.9

$ws.Range("H58").Value = 3334.2
$ws.Range("I58").Value = 1264.7894
$ws.Range("J58").Value = 9887.333000000001
$ws.Range("K58").Value = 1264.7894
$ws.Range("L58").Value = 9887.333000000001
$ws.Range("M58").Value = -1061.7894
$ws.Range("N58").Value = -10293.333

$ws.Range("H74").Value = 101999.3
$ws.Range("J74").Value = 173998.6
$ws.Range("L74").Value = 173998.6
$ws.Range("N74").Value = -175746.6

$ws.Range("H77").Value = 101999.3
$ws.Range("J77").Value = 173998.6
$ws.Range("L77").Value = 521995.8
$ws.Range("N77").Value = -530731.8

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H94").Value = 3430.75
$ws.Range("I94").Value = 2488
$ws.Range("J94").Value = 4750.6
$ws.Range("K94").Value = 2488
$ws.Range("L94").Value = 4750.6
$ws.Range("M94").Value = -2037
$ws.Range("N94").Value = -5652.6

$ws.Range("H132").Value = 3497.48
$ws.Range("I132").Value = 2382.2727
$ws.Range("K132").Value = 7146.8181
$ws.Range("M132").Value = -4616.8181

$ws.Range("H134").Value = 6852.0435
$ws.Range("I134").Value = 4840.5
$ws.Range("J134").Value = 8399.385
$ws.Range("K134").Value = 14521.5
$ws.Range("L134").Value = 25198.155
$ws.Range("M134").Value = -11986.5
$ws.Range("N134").Value = -30268.155

$ws.Range("H135").Value = 67194.14
$ws.Range("J135").Value = 67194.14
$ws.Range("L135").Value = 67194.14
$ws.Range("N135").Value = -77334.14

$ws.Range("H136").Value = 3334.2
$ws.Range("I136").Value = 1264.7894
$ws.Range("J136").Value = 9887.333000000001
$ws.Range("K136").Value = 3794.3682
$ws.Range("L136").Value = 29661.999
$ws.Range("M136").Value = -1244.3682
$ws.Range("N136").Value = -34761.999

$ws.Range("H138").Value = 64930
$ws.Range("J138").Value = 64930
$ws.Range("L138").Value = 64930
$ws.Range("N138").Value = -75210

$ws.Range("H140").Value = 64913.2
$ws.Range("J140").Value = 64913.2
$ws.Range("L140").Value = 64913.2
$ws.Range("N140").Value = -75273.2

$ws.Range("H141").Value = 144994.3
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 144994.3
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 144994.3
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -155354.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 25141.709
$ws.Range("I2").Value = 203
$ws.Range("J2").Value = 35410.59
$ws.Range("K2").Value = 1218
$ws.Range("L2").Value = 212463.54
$ws.Range("M2").Value = -1105
$ws.Range("N2").Value = -212689.54

$ws.Range("H37").Value = 216267.06
$ws.Range("J37").Value = 216267.06
$ws.Range("L37").Value = 648801.1799999999
$ws.Range("N37").Value = -649025.1799999999

$ws.Range("H113").Value = 71429700
$ws.Range("J113").Value = 111112290
$ws.Range("L113").Value = 333336870
$ws.Range("N113").Value = -333341210

$ws.Range("H137").Value = 66636.625
$ws.Range("I137").Value = 1687.3334
$ws.Range("J137").Value = 105606.2
$ws.Range("K137").Value = 5062.0002
$ws.Range("L137").Value = 316818.6
$ws.Range("M137").Value = 37.9997999999996
$ws.Range("N137").Value = -327018.6

$ws.Range("H139").Value = 7986.273
$ws.Range("I139").Value = 2867.75
$ws.Range("J139").Value = 10911.143
$ws.Range("K139").Value = 8603.25
$ws.Range("L139").Value = 32733.429
$ws.Range("M139").Value = -3463.25
$ws.Range("N139").Value = -43013.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2634.6
$ws.Range("J97").Value = 4994
$ws.Range("L97").Value = 4994
$ws.Range("N97").Value = -5986

$ws.Range("H102").Value = 3177501.2
$ws.Range("I102").Value = 3706685.5
$ws.Range("K102").Value = 3706685.5
$ws.Range("M102").Value = -3705063.5

$ws.Range("H122").Value = 16399.25
$ws.Range("I122").Value = 21935.6
$ws.Range("K122").Value = 65806.79999999999
$ws.Range("M122").Value = -63356.79999999999

$ws.Range("H126").Value = 3475.5
$ws.Range("I126").Value = 2184.6365
$ws.Range("J126").Value = 5504
$ws.Range("K126").Value = 6553.9095
$ws.Range("L126").Value = 16512
$ws.Range("M126").Value = -4083.9095
$ws.Range("N126").Value = -21452

$ws.Range("H132").Value = 331535.4
$ws.Range("I132").Value = 377897.22
$ws.Range("K132").Value = 1133691.66
$ws.Range("M132").Value = -1131161.66

$ws.Range("H135").Value = 69995
$ws.Range("J135").Value = 69995
$ws.Range("L135").Value = 69995
$ws.Range("N135").Value = -80135

$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2857.1428
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2857.1428
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2857.1428
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -3081.1428

$ws.Range("H7").Value = 7792.517
$ws.Range("I7").Value = 5116.222
$ws.Range("J7").Value = 12171.909
$ws.Range("K7").Value = 5116.222
$ws.Range("L7").Value = 12171.909
$ws.Range("M7").Value = -5004.222
$ws.Range("N7").Value = -12395.909

$ws.Range("H48").Value = 28000
$ws.Range("I48").Value = 6000
$ws.Range("J48").Value = 50000
$ws.Range("K48").Value = 6000
$ws.Range("L48").Value = 50000
$ws.Range("M48").Value = -5339
$ws.Range("N48").Value = -51322

$ws.Range("H55").Value = 2633287.8
$ws.Range("I55").Value = 6250366
$ws.Range("J55").Value = 2685.3635
$ws.Range("K55").Value = 6250366
$ws.Range("L55").Value = 2685.3635
$ws.Range("M55").Value = -6250193
$ws.Range("N55").Value = -3031.3635

$ws.Range("H126").Value = 7792.517
$ws.Range("I126").Value = 5116.222
$ws.Range("J126").Value = 12171.909
$ws.Range("K126").Value = 15348.666
$ws.Range("L126").Value = 36515.727
$ws.Range("M126").Value = -12878.666
$ws.Range("N126").Value = -41455.727

$ws.Range("H127").Value = 49000
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H130").Value = 64127.453
$ws.Range("J130").Value = 64127.453
$ws.Range("L130").Value = 64127.453
$ws.Range("N130").Value = -74167.45300000001

$ws.Range("H136").Value = 6165.35
$ws.Range("I136").Value = 2753.4666
$ws.Range("K136").Value = 8260.399800000001
$ws.Range("M136").Value = -5710.399800000001

$ws.Range("H141").Value = 69895
$ws.Range("J141").Value = 69895
$ws.Range("L141").Value = 69895
$ws.Range("N141").Value = -80255

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 50000000
$ws.Range("J5").Value = 50000000
$ws.Range("L5").Value = 50000000
$ws.Range("N5").Value = -50000224

$ws.Range("H75").Value = 28748.75
$ws.Range("I75").Value = 26498
$ws.Range("K75").Value = 26498
$ws.Range("M75").Value = -25562

$ws.Range("H78").Value = 28748.75
$ws.Range("I78").Value = 26498
$ws.Range("K78").Value = 79494
$ws.Range("M78").Value = -74814

$ws.Range("H122").Value = 15052.25
$ws.Range("I122").Value = 3200
$ws.Range("J122").Value = 19003
$ws.Range("K122").Value = 9600
$ws.Range("L122").Value = 57009
$ws.Range("M122").Value = -7150
$ws.Range("N122").Value = -61909

$ws.Range("H126").Value = 2366.6128
$ws.Range("I126").Value = 1856
$ws.Range("K126").Value = 5568
$ws.Range("M126").Value = -3098

$ws.Range("H132").Value = 4648.349
$ws.Range("I132").Value = 2126.8914
$ws.Range("K132").Value = 6380.674199999999
$ws.Range("M132").Value = -3850.674199999999

$ws.Range("H135").Value = 60830.668
$ws.Range("J135").Value = 60830.668
$ws.Range("L135").Value = 60830.668
$ws.Range("N135").Value = -70970.66800000001

$ws.Range("H136").Value = 2335.457
$ws.Range("I136").Value = 1625.1724
$ws.Range("J136").Value = 5768.5
$ws.Range("K136").Value = 4875.5172
$ws.Range("L136").Value = 17305.5
$ws.Range("M136").Value = -2325.5172
$ws.Range("N136").Value = -22405.5

$ws.Range("H137").Value = 69995
$ws.Range("J137").Value = 69995
$ws.Range("L137").Value = 69995
$ws.Range("N137").Value = -80195

$ws.Range("H139").Value = 68783.5
$ws.Range("J139").Value = 70355
$ws.Range("L139").Value = 70355
$ws.Range("N139").Value = -80635

$ws.Range("H141").Value = 104169.9
$ws.Range("J141").Value = 104169.9
$ws.Range("L141").Value = 104169.9
$ws.Range("N141").Value = -114529.9
